$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 78 - Yi Large
$ws.Range("A78").Value = "Yi Large"
$ws.Range("B78").Value = 45425
$ws.Range("C78").Value = "01.AI"
$ws.Range("D78").Value = "100*"
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = 16
$ws.Range("G78").Value = 2.8
$ws.Range("H78").Value = 2.8
$ws.Range("K78").Value = "Proprietary"
$ws.Range("L78").Value = $true
$ws.Range("M78").Value = $true
$ws.Range("P78").Value = 83.8
$ws.Range("R78").Value = 82.3
$ws.Range("S78").Value = 62.4

# Row 79 - Yi Medium
$ws.Range("A79").Value = "Yi Medium"
$ws.Range("B79").Value = 45425
$ws.Range("C79").Value = "01.AI"
$ws.Range("D79").Value = "10*"
$ws.Range("E79").Value = 16
$ws.Range("F79").Value = 16
$ws.Range("G79").Value = 0.35
$ws.Range("H79").Value = 0.35
$ws.Range("K79").Value = "Proprietary"
$ws.Range("P79").Value = 76.8
$ws.Range("R79").Value = 75.2
$ws.Range("S79").Value = 50.1

# Row 80 - GLM-4
# B80 is a brand-new row/cell (rows 80-91 did not exist before), so copy the
# date number format from an existing date cell (e.g. B78) before writing the value.
$ws.Range("B78").Copy()
$ws.Range("B80").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A80").Value = "GLM-4"
$ws.Range("B80").Value = 45307
$ws.Range("C80").Value = "Zhipu AI"
$ws.Range("D80").Value = "130*"
$ws.Range("E80").Value = 128
$ws.Range("F80").Value = 8.192
$ws.Range("G80").Value = 14.1
$ws.Range("H80").Value = 14.1
$ws.Range("K80").Value = "Proprietary"
$ws.Range("L80").Value = $true
$ws.Range("M80").Value = $true
$ws.Range("P80").Value = 81.5
$ws.Range("R80").Value = 72
$ws.Range("S80").Value = 47.9

# Row 81 - GLM-4V
$ws.Range("B78").Copy()
$ws.Range("B81").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A81").Value = "GLM-4V"
$ws.Range("B81").Value = 45307
$ws.Range("C81").Value = "Zhipu AI"
$ws.Range("D81").Value = "130*"
$ws.Range("E81").Value = 2
$ws.Range("F81").Value = 2
$ws.Range("G81").Value = 14.1
$ws.Range("H81").Value = 14.1
$ws.Range("I81").Value = 14.8
$ws.Range("K81").Value = "Proprietary"
$ws.Range("L81").Value = $true
$ws.Range("M81").Value = $true
$ws.Range("N81").Value = $true

# Update frozen pane / selection view state
$ws.Application.ActiveWindow.ScrollRow = 68
$ws.Range("I84").Select()
